$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before the current row 859 (old rows 859:873 shift down to 865:879)
$ws.Rows("859:864").Insert()

# Constant values shared by every data row in this sheet
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112037
$categoria = "Cebollín"
$variedad = "Sin especificar"
$unidad = "`$/paquete 36 unidades"
$kgUnidades = 36
$clasificacion = "Hortaliza"

# New rows: a new weekly report (Fecha = 2021-09-09 -> serial 44448)
$newRows = @(
  @{ Row=859; Fecha=44448; Calidad="Extra";   Volumen=810; PrecioMin=3000; PrecioMax=3300; PrecioProm=3163; Origen="Provincia de Chacabuco"; PrecioKg=88 },
  @{ Row=860; Fecha=44448; Calidad="Extra";   Volumen=720; PrecioMin=3000; PrecioMax=3300; PrecioProm=3175; Origen="Región Metropolitana";   PrecioKg=88 },
  @{ Row=861; Fecha=44448; Calidad="Primera"; Volumen=820; PrecioMin=2200; PrecioMax=2500; PrecioProm=2390; Origen="Provincia de Chacabuco"; PrecioKg=66 },
  @{ Row=862; Fecha=44448; Calidad="Primera"; Volumen=860; PrecioMin=2300; PrecioMax=2500; PrecioProm=2409; Origen="Región Metropolitana";   PrecioKg=67 },
  @{ Row=863; Fecha=44448; Calidad="Segunda"; Volumen=250; PrecioMin=1800; PrecioMax=1800; PrecioProm=1800; Origen="Provincia de Chacabuco"; PrecioKg=50 },
  @{ Row=864; Fecha=44448; Calidad="Segunda"; Volumen=330; PrecioMin=1900; PrecioMax=1900; PrecioProm=1900; Origen="Región Metropolitana";   PrecioKg=53 }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $mercadoId
  $ws.Cells.Item($row, 2).Value = $mercado
  $ws.Cells.Item($row, 3).Value = $region
  $ws.Cells.Item($row, 4).Value = $r.Fecha
  $ws.Cells.Item($row, 5).Value = $codreg
  $ws.Cells.Item($row, 6).Value = $categoriaId
  $ws.Cells.Item($row, 7).Value = $categoria
  $ws.Cells.Item($row, 8).Value = $variedad
  $ws.Cells.Item($row, 9).Value = $r.Calidad
  $ws.Cells.Item($row, 10).Value = $r.Volumen
  $ws.Cells.Item($row, 11).Value = $r.PrecioMin
  $ws.Cells.Item($row, 12).Value = $r.PrecioMax
  $ws.Cells.Item($row, 13).Value = $r.PrecioProm
  $ws.Cells.Item($row, 14).Value = $unidad
  $ws.Cells.Item($row, 15).Value = $r.Origen
  $ws.Cells.Item($row, 16).Value = $r.PrecioKg
  $ws.Cells.Item($row, 17).Value = $kgUnidades
  $ws.Cells.Item($row, 18).Value = $clasificacion
}
